$wb = $excel.ActiveWorkbook
$wsDone = $wb.Worksheets.Item("done")
$wsShort = $wb.Worksheets.Item("short term")

# --- "short term" sheet: remove the row that held the old, short version of
# item 67 ("tenter de fixer le f a 0 ..."); the fuller rewrite of that task
# is being moved to the "done" sheet below.
$wsShort.Rows(19).Delete()

# --- "done" sheet: append a new task row with the reworked / expanded
# wording for item 67, copying the formatting from the row above it (A48)
# so the new row matches the look of the rest of the list.
$wsDone.Range("A48").Copy()
$wsDone.Range("A49").PasteSpecial(-4122)  # xlPasteFormats

$newTaskText = "67. tenter de fixer le f à 0 pour les modèles Gauss-probit et log-Gauss-probit, sans réajuster ou en réajustant à partir des mêmes valeurs, et garder le probit ou log-probit s'il est meilleur en AIC: réajustement systématique avec f = 0  du modèle si GP5p ou lGP5p et retenue du modèle simplifié sur critère d'information meilleur"

$wsDone.Range("A49").Value = $newTaskText
$wsDone.Range("B49").Value = "ML"

$wsDone.Range("A49").Select()
